$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("filas")
$ws2 = $wb.Worksheets.Item("columnas")

for ($r = 8; $r -le 19; $r++) {
    $ws1.Cells.Item($r, 6).Value = 2
    $ws1.Cells.Item($r, 7).Value = "Utilización"
}

$ws1.Range("F1").Select() | Out-Null
$ws2.Activate() | Out-Null
